$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Comment" (column E) status text for each row. New distinct
# phrases are introduced group-by-group (not strictly row order) so that the
# shared-string table ends up with the same new entries in the same order as
# the target workbook.

# 'Alias Changed for Network'
foreach ($row in @(2, 3, 4, 5, 6, 7, 8)) {
    $ws.Cells.Item($row, 5).Value = 'Alias Changed for Network'
}

# 'Name of Add-On Package Changed'
foreach ($row in @(33, 34, 35, 39, 46, 52, 57, 60, 61)) {
    $ws.Cells.Item($row, 5).Value = 'Name of Add-On Package Changed'
}

# 'New Network Added to Database in Aug 2020'
foreach ($row in @(64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76)) {
    $ws.Cells.Item($row, 5).Value = 'New Network Added to Database in Aug 2020'
}

# 'Network Added to Add-On Package in Aug 2020'
foreach ($row in @(32, 42, 43, 44, 47, 50, 53, 54, 59, 63, 78, 86)) {
    $ws.Cells.Item($row, 5).Value = 'Network Added to Add-On Package in Aug 2020'
}

# 'Network Removed from Add-On Package in Aug 2020'
foreach ($row in @(62, 79, 81, 83, 84, 85, 91, 92, 93)) {
    $ws.Cells.Item($row, 5).Value = 'Network Removed from Add-On Package in Aug 2020'
}

# 'Network Moved from Base Service to Add-On Package in Aug 2020'
foreach ($row in @(17, 23, 58, 88, 94)) {
    $ws.Cells.Item($row, 5).Value = 'Network Moved from Base Service to Add-On Package in Aug 2020'
}

# Update the saved selection on the sheet (was A23) to B2.
$ws.Range("B2").Select()
